$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neurology")

# Delete the last three data rows (30-32), removing the trailing 17:01 scan
# entries for students 212572 / 213007 / 213006. Rows below shift up, so the
# sheet's dimension shrinks from A1:F32 to A1:F29.
$ws.Rows.Item(30).Resize(3).Delete()

# Rename the worksheet/tab from "Neurology" to "Session"
$ws.Name = "Session"
